$wb = $excel.ActiveWorkbook

# ALC row 32 (G32=5484)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1206.6666
$ws.Range("I32").Value = 775
$ws.Range("J32").Value = 1422.5
$ws.Range("K32").Value = 775
$ws.Range("L32").Value = 1422.5
$ws.Range("M32").Value = -449
$ws.Range("N32").Value = -2074.5

# ALC row 64 (G64=5506)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 6857.143
$ws.Range("I64").Value = 8140
$ws.Range("K64").Value = 8140
$ws.Range("M64").Value = -7892

# ALC row 67 (G67=5506)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 6857.143
$ws.Range("I67").Value = 8140
$ws.Range("K67").Value = 8140
$ws.Range("M67").Value = -7282

# ALC row 92 (G92=19901)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 2694.3
$ws.Range("I92").Value = 3036.3125
$ws.Range("J92").Value = 1326.25
$ws.Range("K92").Value = 3036.3125
$ws.Range("L92").Value = 1326.25
$ws.Range("M92").Value = -1788.3125
$ws.Range("N92").Value = -3822.25

# ALC row 106 (G106=19903)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2184.1667
$ws.Range("I106").Value = 1968.3334
$ws.Range("J106").Value = 2400
$ws.Range("K106").Value = 1968.3334
$ws.Range("L106").Value = 2400
$ws.Range("M106").Value = -1337.3334
$ws.Range("N106").Value = -3662

# ALC row 112 (G112=27960)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 4132.143
$ws.Range("I112").Value = 2200
$ws.Range("J112").Value = 4364
$ws.Range("K112").Value = 6600
$ws.Range("L112").Value = 13092
$ws.Range("M112").Value = -5492
$ws.Range("N112").Value = -15308

# ALC row 136 (G136=42164)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 82300
$ws.Range("J136").Value = 84950
$ws.Range("L136").Value = 84950
$ws.Range("N136").Value = -95150

# ARM row 32 (G32=44147)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20616.707
$ws.Range("I32").Value = 15991.642
$ws.Range("J32").Value = 59351.625
$ws.Range("K32").Value = 15991.642
$ws.Range("L32").Value = 59351.625
$ws.Range("M32").Value = -15704.642
$ws.Range("N32").Value = -59925.625

# ARM row 61 (G61=43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1118.625
$ws.Range("I61").Value = 1049.1482
$ws.Range("J61").Value = 1493.8
$ws.Range("K61").Value = 1049.1482
$ws.Range("L61").Value = 1493.8
$ws.Range("M61").Value = -837.1482000000001
$ws.Range("N61").Value = -1917.8

# ARM row 63 (G63=12528)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2753.257
$ws.Range("I63").Value = 2172.7407
$ws.Range("J63").Value = 4712.5
$ws.Range("K63").Value = 2172.7407
$ws.Range("L63").Value = 4712.5
$ws.Range("M63").Value = -1486.7407
$ws.Range("N63").Value = -6084.5

# ARM row 66 (G66=12528)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2753.257
$ws.Range("I66").Value = 2172.7407
$ws.Range("J66").Value = 4712.5
$ws.Range("K66").Value = 10863.7035
$ws.Range("L66").Value = 23562.5
$ws.Range("M66").Value = -7431.7035
$ws.Range("N66").Value = -30426.5

# ARM row 97 (G97=19941)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2250.5789
$ws.Range("I97").Value = 2603.2222
$ws.Range("J97").Value = 1385
$ws.Range("K97").Value = 2603.2222
$ws.Range("L97").Value = 1385
$ws.Range("M97").Value = -2107.2222
$ws.Range("N97").Value = -2377

# ARM row 136 (G136=43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1118.625
$ws.Range("I136").Value = 1049.1482
$ws.Range("J136").Value = 1493.8
$ws.Range("K136").Value = 3147.4446
$ws.Range("L136").Value = 4481.4
$ws.Range("M136").Value = -597.4446000000003
$ws.Range("N136").Value = -9581.4

# BSM row 20 (G20=14149)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 35741092
$ws.Range("I20").Value = 44443.93
$ws.Range("K20").Value = 44443.93
$ws.Range("M20").Value = -44196.93

# BSM row 94 (G94=19939)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 615.6667
$ws.Range("I94").Value = 501.5
$ws.Range("J94").Value = 1118
$ws.Range("K94").Value = 501.5
$ws.Range("L94").Value = 1118
$ws.Range("M94").Value = -50.5
$ws.Range("N94").Value = -2020

# BSM row 105 (G105=19947)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2499.3462
$ws.Range("I105").Value = 2541.279
$ws.Range("J105").Value = 2299
$ws.Range("K105").Value = 2541.279
$ws.Range("L105").Value = 2299
$ws.Range("M105").Value = -794.279
$ws.Range("N105").Value = -5793

# BSM row 107 (G107=27706)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 8835.235000000001
$ws.Range("I107").Value = 1290.9
$ws.Range("J107").Value = 19612.857
$ws.Range("K107").Value = 1290.9
$ws.Range("L107").Value = 19612.857
$ws.Range("M107").Value = 629.0999999999999
$ws.Range("N107").Value = -23452.857

# BSM row 141 (G141=43278)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H141").Value = 61565.25
$ws.Range("J141").Value = 61565.25
$ws.Range("L141").Value = 61565.25
$ws.Range("N141").Value = -71925.25

# CRP row 31 (G31=44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3281
$ws.Range("I31").Value = 2826.862
$ws.Range("J31").Value = 5162.4287
$ws.Range("K31").Value = 2826.862
$ws.Range("L31").Value = 5162.4287
$ws.Range("M31").Value = -2531.862
$ws.Range("N31").Value = -5752.4287

# CRP row 34 (G34=44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3281
$ws.Range("I34").Value = 2826.862
$ws.Range("J34").Value = 5162.4287
$ws.Range("K34").Value = 2826.862
$ws.Range("L34").Value = 5162.4287
$ws.Range("M34").Value = -2624.862
$ws.Range("N34").Value = -5566.4287

# CRP row 105 (G105=19928)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 949.75
$ws.Range("I105").Value = 949.75
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 949.75
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 797.25
$ws.Range("N105").ClearContents()

# CRP row 112 (G112=25796)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H112").Value = 25222.5
$ws.Range("J112").Value = 25222.5
$ws.Range("L112").Value = 25222.5
$ws.Range("N112").Value = -28176.5

# CUL row 131 (G131=36060)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 119884.12
$ws.Range("I131").Value = 336933.34
$ws.Range("J131").Value = 73373.57000000001
$ws.Range("K131").Value = 1010800.02
$ws.Range("L131").Value = 220120.71
$ws.Range("M131").Value = -1005760.02
$ws.Range("N131").Value = -230200.71

# GSM row 11 (G11=4422)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 18338334
$ws.Range("I11").Value = 22000000
$ws.Range("J11").Value = 30000
$ws.Range("K11").Value = 22000000
$ws.Range("L11").Value = 30000
$ws.Range("M11").Value = -21999861
$ws.Range("N11").Value = -30278

# GSM row 122 (G122=36182)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 488750.94
$ws.Range("I122").Value = 775189.1
$ws.Range("J122").Value = 1806
$ws.Range("K122").Value = 2325567.3
$ws.Range("L122").Value = 5418
$ws.Range("M122").Value = -2323117.3
$ws.Range("N122").Value = -10318

# GSM row 126 (G126=36184)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3783.1333
$ws.Range("I126").Value = 2672.8462
$ws.Range("J126").Value = 11000
$ws.Range("K126").Value = 8018.5386
$ws.Range("L126").Value = 33000
$ws.Range("M126").Value = -5548.5386
$ws.Range("N126").Value = -37940

# LTW row 40 (G40=36248)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1685984
$ws.Range("I40").Value = 3369134.8
$ws.Range("J40").Value = 2833.3333
$ws.Range("K40").Value = 3369134.8
$ws.Range("L40").Value = 2833.3333
$ws.Range("M40").Value = -3368998.8
$ws.Range("N40").Value = -3105.3333

# LTW row 82 (G82=12565)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3100.1765
$ws.Range("I82").Value = 2999.9092
$ws.Range("J82").Value = 3284
$ws.Range("K82").Value = 2999.9092
$ws.Range("L82").Value = 3284
$ws.Range("M82").Value = -2638.9092
$ws.Range("N82").Value = -4006

# LTW row 85 (G85=12565)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 3100.1765
$ws.Range("I85").Value = 2999.9092
$ws.Range("J85").Value = 3284
$ws.Range("K85").Value = 2999.9092
$ws.Range("L85").Value = 3284
$ws.Range("M85").Value = -1751.9092
$ws.Range("N85").Value = -5780

# LTW row 122 (G122=36247)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2910.1936
$ws.Range("I122").Value = 2695.0557
$ws.Range("J122").Value = 3208.077
$ws.Range("K122").Value = 8085.1671
$ws.Range("L122").Value = 9624.231
$ws.Range("M122").Value = -5635.1671
$ws.Range("N122").Value = -14524.231

# LTW row 132 (G132=44058)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3746.7932
$ws.Range("I132").Value = 2631.111
$ws.Range("J132").Value = 5572.4546
$ws.Range("K132").Value = 7893.333
$ws.Range("L132").Value = 16717.3638
$ws.Range("M132").Value = -5363.333
$ws.Range("N132").Value = -21777.3638

# WVR row 96 (G96=19977)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3878.6316
$ws.Range("I96").Value = 2814.2856
$ws.Range("J96").Value = 4499.5
$ws.Range("K96").Value = 2814.2856
$ws.Range("L96").Value = 4499.5
$ws.Range("M96").Value = -1441.2856
$ws.Range("N96").Value = -7245.5

# WVR row 126 (G126=36210)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 931.5454999999999
$ws.Range("I126").Value = 931.5454999999999
$ws.Range("K126").Value = 2794.6365
$ws.Range("M126").Value = -324.6364999999996

# WVR row 132 (G132=44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3202.6223
$ws.Range("I132").Value = 994.4167
$ws.Range("J132").Value = 12035.444
$ws.Range("K132").Value = 2983.2501
$ws.Range("L132").Value = 36106.33199999999
$ws.Range("M132").Value = -453.2501000000002
$ws.Range("N132").Value = -41166.33199999999

